# Minor wording edits on slide 1's "The Problem:" body text box.
# Splits two paragraphs into multiple runs (same rPr) reflecting small
# copy-edits:
#   - "non-experts technical"        -> "non-experts, technical"
#   - "As an expert is can be"       -> "As an expert, it can be"
#   - "have successful removed"      -> "have successfully removed"

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# --- Paragraph 1: "When communicating with non-experts technical language or jargon creates barriers for understanding." ---
$para1 = $tr.Paragraphs(1, 1)

# Work right-to-left so already-computed character offsets (from the
# original text) stay valid as each replacement shifts the text after it.

# "technical language or jargon creates barriers for understanding." (chars 37-100, len 64)
$para1.Characters(37, 64).Text = "technical language or jargon creates barriers for understanding."

# "experts " -> "experts, " (chars 29-36, len 8)
$para1.Characters(29, 8).Text = "experts, "

# --- Paragraph 3: "As an expert is can be a challenge to identify if you have successful removed jargon from a document." ---
$para3 = $tr.Paragraphs(3, 1)

# "removed jargon from a document." (chars 71-101, len 31)
$para3.Characters(71, 31).Text = "removed jargon from a document."

# "successful " -> "successfully " (chars 60-70, len 11)
$para3.Characters(60, 11).Text = "successfully "

# "can be a challenge to identify if you have " (chars 17-59, len 43)
$para3.Characters(17, 43).Text = "can be a challenge to identify if you have "

# "expert is " -> "expert, it " (chars 7-16, len 10)
$para3.Characters(7, 10).Text = "expert, it "
